$d = $word.ActiveDocument

# Locate the stray paragraph that contains only the italic "2 Pedro" run
# (sitting right after the "2PE" Heading2 paragraph) and remove the whole
# paragraph, including its paragraph mark.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "2 Pedro" -and $p.Range.Font.Italic -eq -1) {
        $p.Range.Delete()
        break
    }
}
